$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 5
$ws.Range("H5").Value = 154.75
$ws.Range("I5").Value = 173
$ws.Range("K5").Value = 173
$ws.Range("M5").Value = -58

# Row 32
$ws.Range("H32").Value = 20244.75
$ws.Range("I32").Value = 18979
$ws.Range("K32").Value = 18979
$ws.Range("M32").Value = -18653

# Row 40
$ws.Range("H40").Value = 4274.48
$ws.Range("J40").Value = 4313.7085
$ws.Range("L40").Value = 4313.7085
$ws.Range("N40").Value = -4663.7085

# Row 43
$ws.Range("H43").Value = 5511
$ws.Range("I43").Value = 4931
$ws.Range("J43").Value = 5859
$ws.Range("K43").Value = 4931
$ws.Range("L43").Value = 5859
$ws.Range("M43").Value = -4862
$ws.Range("N43").Value = -5997

# Row 51
$ws.Range("H51").Value = 14016.777
$ws.Range("J51").Value = 3899
$ws.Range("L51").Value = 3899
$ws.Range("N51").Value = -4867

# Row 52
$ws.Range("H52").Value = 364.3889
$ws.Range("I52").Value = 129.5
$ws.Range("J52").Value = 393.75
$ws.Range("K52").Value = 388.5
$ws.Range("L52").Value = 1181.25
$ws.Range("M52").Value = -228.5
$ws.Range("N52").Value = -1501.25

# Row 86
$ws.Range("H86").Value = 2454.3809
$ws.Range("I86").Value = 2330.8
$ws.Range("J86").Value = 2493
$ws.Range("K86").Value = 2330.8
$ws.Range("L86").Value = 2493
$ws.Range("M86").Value = -1207.8
$ws.Range("N86").Value = -4739

# Row 87
$ws.Range("H87").Value = 27071.4
$ws.Range("J87").Value = 27071.4
$ws.Range("L87").Value = 27071.4
$ws.Range("N87").Value = -29567.4

# Row 89
$ws.Range("H89").Value = 2454.3809
$ws.Range("I89").Value = 2330.8
$ws.Range("J89").Value = 2493
$ws.Range("K89").Value = 11654
$ws.Range("L89").Value = 12465
$ws.Range("M89").Value = -6038
$ws.Range("N89").Value = -23697

# Row 90
$ws.Range("H90").Value = 27071.4
$ws.Range("J90").Value = 27071.4
$ws.Range("L90").Value = 81214.20000000001
$ws.Range("N90").Value = -93694.20000000001

# Row 96
$ws.Range("H96").Value = 3060.4
$ws.Range("I96").Value = 1784
$ws.Range("J96").Value = 4975
$ws.Range("K96").Value = 5352
$ws.Range("L96").Value = 14925
$ws.Range("M96").Value = -3979
$ws.Range("N96").Value = -17671

# Row 132
$ws.Range("H132").Value = 1498
$ws.Range("I132").Value = 1498
$ws.Range("K132").Value = 4494
$ws.Range("M132").Value = -1964

# Row 138
$ws.Range("H138").Value = 3524.4614
$ws.Range("J138").Value = 3964.1516
$ws.Range("L138").Value = 11892.4548
$ws.Range("N138").Value = -22172.4548

# Row 141
$ws.Range("H141").Value = 77424.53999999999
$ws.Range("I141").Value = 83459.914
$ws.Range("K141").Value = 250379.742
$ws.Range("M141").Value = -245199.742


$ws = $wb.Worksheets.Item("ARM")
# Row 34
$ws.Range("H34").Value = 125000
$ws.Range("I34").Value = 125000
$ws.Range("K34").Value = 125000
$ws.Range("M34").Value = -124729

# Row 88
$ws.Range("H88").Value = 1646.7693
$ws.Range("J88").Value = 2601.5
$ws.Range("L88").Value = 2601.5
$ws.Range("N88").Value = -3413.5

# Row 91
$ws.Range("H91").Value = 1646.7693
$ws.Range("J91").Value = 2601.5
$ws.Range("L91").Value = 2601.5
$ws.Range("N91").Value = -5409.5

# Row 97
$ws.Range("H97").Value = 4802
$ws.Range("I97").Value = 5189.9375
$ws.Range("K97").Value = 5189.9375
$ws.Range("M97").Value = -4693.9375


$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 65271.668
$ws.Range("I86").Value = 204658.6
$ws.Range("K86").Value = 204658.6
$ws.Range("M86").Value = -203535.6

# Row 89
$ws.Range("H89").Value = 65271.668
$ws.Range("I89").Value = 204658.6
$ws.Range("K89").Value = 1023293
$ws.Range("M89").Value = -1017677

# Row 99
$ws.Range("H99").Value = 3000
$ws.Range("I99").Value = 3000
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 3000
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1502
$ws.Range("N99").ClearContents() | Out-Null

# Row 105
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").ClearContents() | Out-Null

# Row 135
$ws.Range("H135").Value = 49999.5
$ws.Range("J135").Value = 49999.5
$ws.Range("L135").Value = 49999.5
$ws.Range("N135").Value = -60139.5


$ws = $wb.Worksheets.Item("CRP")
# Row 25
$ws.Range("H25").Value = 6495
$ws.Range("I25").Value = 6495
$ws.Range("K25").Value = 6495
$ws.Range("M25").Value = -6321

# Row 68
$ws.Range("H68").Value = 23000
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 23000
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 23000
$ws.Range("M68").ClearContents() | Out-Null
$ws.Range("N68").Value = -24498

# Row 71
$ws.Range("H71").Value = 23000
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 23000
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 69000
$ws.Range("M71").ClearContents() | Out-Null
$ws.Range("N71").Value = -76488

# Row 135
$ws.Range("H135").Value = 80779.5
$ws.Range("J135").Value = 80779.5
$ws.Range("L135").Value = 80779.5
$ws.Range("N135").Value = -90919.5


$ws = $wb.Worksheets.Item("CUL")
# Row 11
$ws.Range("H11").Value = 1076.7916
$ws.Range("I11").Value = 1115.7273
$ws.Range("K11").Value = 3347.1819
$ws.Range("M11").Value = -3207.1819

# Row 12
$ws.Range("H12").Value = 252.025
$ws.Range("I12").Value = 110.64286
$ws.Range("J12").Value = 581.9167
$ws.Range("K12").Value = 331.92858
$ws.Range("L12").Value = 1745.7501
$ws.Range("M12").Value = -158.92858
$ws.Range("N12").Value = -2091.7501

# Row 23
$ws.Range("H23").Value = 189.90909
$ws.Range("I23").Value = 197.5
$ws.Range("J23").Value = 185.57143
$ws.Range("K23").Value = 592.5
$ws.Range("L23").Value = 556.71429
$ws.Range("M23").Value = -357.5
$ws.Range("N23").Value = -1026.71429


$ws = $wb.Worksheets.Item("GSM")
# Row 10
$ws.Range("H10").Value = 5030275
$ws.Range("J10").Value = 10550
$ws.Range("L10").Value = 10550
$ws.Range("N10").Value = -10888

# Row 11
$ws.Range("H11").Value = 2232799.8
$ws.Range("J11").Value = 6532.8335
$ws.Range("L11").Value = 6532.8335
$ws.Range("N11").Value = -6810.8335

# Row 18
$ws.Range("H18").Value = 15999.5
$ws.Range("J18").Value = 15999.5
$ws.Range("L18").Value = 15999.5
$ws.Range("N18").Value = -16585.5

# Row 80
$ws.Range("H80").Value = 4009.4092
$ws.Range("J80").Value = 4820
$ws.Range("L80").Value = 4820
$ws.Range("N80").Value = -6816

# Row 83
$ws.Range("H83").Value = 4009.4092
$ws.Range("J83").Value = 4820
$ws.Range("L83").Value = 24100
$ws.Range("N83").Value = -34084

# Row 132
$ws.Range("H132").Value = 1453.5
$ws.Range("I132").Value = 1323.3846
$ws.Range("K132").Value = 3970.1538
$ws.Range("M132").Value = -1440.1538

# Row 135
$ws.Range("H135").Value = 78223.8
$ws.Range("J135").Value = 78223.8
$ws.Range("L135").Value = 78223.8
$ws.Range("N135").Value = -88363.8


$ws = $wb.Worksheets.Item("LTW")
# Row 20
$ws.Range("H20").Value = 444050
$ws.Range("I20").Value = 887000
$ws.Range("J20").Value = 1100
$ws.Range("K20").Value = 887000
$ws.Range("L20").Value = 1100
$ws.Range("M20").Value = -886774
$ws.Range("N20").Value = -1552

# Row 46
$ws.Range("H46").Value = 3386.2
$ws.Range("J46").Value = 3784.2942
$ws.Range("L46").Value = 3784.2942
$ws.Range("N46").Value = -4160.2942

# Row 55
$ws.Range("H55").Value = 257.07693
$ws.Range("I55").Value = 278.22223
$ws.Range("K55").Value = 278.22223
$ws.Range("M55").Value = -105.22223


$ws = $wb.Worksheets.Item("WVR")
# Row 37
$ws.Range("H37").Value = 19475.25
$ws.Range("J37").Value = 15967
$ws.Range("L37").Value = 15967
$ws.Range("N37").Value = -16373

# Row 135
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents() | Out-Null

